$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.633893966674805
$ws.Range("B1").Value = 3.23035740852356
$ws.Range("C1").Value = 4.580103874206543
$ws.Range("D1").Value = 1.378569841384888
$ws.Range("E1").Value = 0.8034419417381287
